$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 33.575
$ws.Range("C2").Value = 20.815
$ws.Range("D2").Value = 44

$ws.Range("B3").Value = 122.035
$ws.Range("C3").Value = 63.27500000000001
$ws.Range("D3").Value = 67
